# "Correzioni ad organizzazione del lavoro"
#
# Applies the content edits made to slide 13 ("Norme per l'organizzazione del
# lavoro") and slide 15 ("Progettazione di dettaglio e Codifica") of the
# presentation: re-worded headings/labels and a couple of accompanying
# textbox resizes/repositions.

function EmuToPt($emu) {
    # PowerPoint's Shape.Left/Top/Width/Height are expressed in points while
    # the underlying OOXML stores EMUs (1 pt = 12700 EMU). The host truncates
    # towards zero when converting back to EMU, so nudge by half an EMU to
    # land on the exact integer instead of the one below it.
    return ([double]$emu + 0.5) / 12700.0
}

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 13
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)

# "Per permettere di tenere traccia del lavoro del team" -> "Tracciamento del
# lavoro del team"
$sh = $s13.Shapes.Item(2)
$sh.TextFrame.TextRange.Text = "Tracciamento del lavoro del team"

# "Per versionamento " -> "Versionamento " (drop the leading "Per " run,
# capitalise the remaining word in place so it stays a single run).
$sh = $s13.Shapes.Item(4)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 4).Text = ""
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 13).Text = "Versionamento"

# "Norme per l'organizzazione del lavoro" -> "Norme e strumenti" + reflow box
$sh = $s13.Shapes.Item(7)
$tr = $sh.TextFrame.TextRange
$tr.Characters(7, 31).Text = "e strumenti"
$sh.Left   = EmuToPt 1573008
$sh.Top    = EmuToPt 274806
$sh.Width  = EmuToPt 6250615
$sh.Height = EmuToPt 584775

# Email textbox just shifts right/up a bit and narrows
$sh = $s13.Shapes.Item(16)
$sh.Left   = EmuToPt 5965391
$sh.Top    = EmuToPt 3714586
$sh.Width  = EmuToPt 4768647
$sh.Height = EmuToPt 400110

# ---------------------------------------------------------------------------
# Slide 15
# ---------------------------------------------------------------------------
$s15 = $p.Slides.Item(15)

# "Progettazione di dettaglio e Codifica" heading: bump the first word's
# size to match the rest of the title box, and grow the box a touch.
$sh = $s15.Shapes.Item(5)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 13).Font.Size = 32
$sh.Height = EmuToPt 1077218
